$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Locator / text corrections (shared-string swaps) ---
$ws.Range("C24").Value = "sDB104-2311002"
$ws.Range("C19").Value = "27 Nov 2023 - 03 Dec 2023"
$ws.Range("D24").Value = "cDB104-2311002"
$ws.Range("N24").Value = "Processing"

# --- Order date (C18) ---
$ws.Range("C18").Value = 45250.0

# --- New column T (20): widen to match S/R (raw width 18.75 chars) ---
$ws.Columns.Item(20).ColumnWidth = 18

# --- New "Estimated Inbound Date" header cell T21, merged T21:T22 ---
$ws.Range("T21").Value = "Estimated Inbound Date"
$ws.Range("T21:T22").Merge()
# Merge() resets borders on the merged cells; restore the look of the
# neighbouring R21:S22 merge block by re-pasting its formatting.
$ws.Range("S21").Copy()
$ws.Range("T21").PasteSpecial($xlPasteFormats)
$ws.Range("S22").Copy()
$ws.Range("T22").PasteSpecial($xlPasteFormats)

# --- Row 23 date values: O23 updated, R23/S23 shift values, T23 gets old R23 value ---
$ws.Range("R23").Copy()
$ws.Range("T23").PasteSpecial($xlPasteFormats)

$ws.Range("O23").Value = 45252.0
$ws.Range("R23").Value = 45280.0
$ws.Range("S23").Value = 45311.0
$ws.Range("T23").Value = 45261.0

# --- Row 24 values ---
$ws.Range("S24").Copy()
$ws.Range("T24").PasteSpecial($xlPasteFormats)

$ws.Range("O24").Value = 0.0
$ws.Range("Q24").Value = 0.0
$ws.Range("S24").Value = 660.0
$ws.Range("T24").Value = 660.0

# --- Row 25 values ---
$ws.Range("R25").Copy()
$ws.Range("T25").PasteSpecial($xlPasteFormats)

$ws.Range("O25").Value = 0.0
$ws.Range("Q25").Value = 0.0
$ws.Range("T25").Value = 660.0

# --- Row 26 values ---
$ws.Range("R26").Copy()
$ws.Range("T26").PasteSpecial($xlPasteFormats)

$ws.Range("O26").Value = 0.0
$ws.Range("Q26").Value = 0.0
$ws.Range("T26").Value = 660.0

# --- Move "AUTHORIZATION" / "Electronically Approved" text from column Q to R (rows 29-30) ---
$val29 = $ws.Range("Q29").Value
$ws.Range("Q29").Copy()
$ws.Range("R29").PasteSpecial($xlPasteFormats)
$ws.Range("R29").Value = $val29
$ws.Range("Q29").Clear()

$val30 = $ws.Range("Q30").Value
$ws.Range("Q30").Copy()
$ws.Range("R30").PasteSpecial($xlPasteFormats)
$ws.Range("R30").Value = $val30
$ws.Range("Q30").Clear()
